$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2021" column (J) is being appended, mirroring the existing per-year
# columns D:I. Each new cell needs to pick up the formatting already used by
# its row, so for rows whose neighbouring (I-column) formatting differs from
# the plain column default, first clone the format from a same-styled cell
# and only then write the value - exactly like the rest of the sheet does.

function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# Divider row (thin bottom border across the whole row).
Copy-Format "I3" "J3"

# Year header row (2015 ... 2020 -> add 2021).
Copy-Format "I4" "J4"
$ws.Range("J4").Value = 2021

# "Median value" row - picks up its own row default style automatically.
$ws.Range("J5").Value = 5356.3

# "Urbanisation" ratio row - picks up its own row default style automatically.
$ws.Range("J6").Value = 9.5

# Plain data rows - these inherit the column default format automatically.
$ws.Range("J8").Value = 7.9
$ws.Range("J9").Value = 10.5
$ws.Range("J11").Value = 9.6
$ws.Range("J12").Value = 9.4
$ws.Range("J14").Value = 14.8
$ws.Range("J15").Value = 9.1
$ws.Range("J16").Value = 9.5
$ws.Range("J17").Value = 5.9

# Region rows (19-26) use the number-formatted style carried by column I,
# minus its border, which matches the existing "28" style elsewhere in the
# sheet - clone it from a cell that already has that exact style.
Copy-Format "I8" "J19"
$ws.Range("J19").Value = 12.434613462352335
Copy-Format "I8" "J20"
$ws.Range("J20").Value = 16.80050595536094
Copy-Format "I8" "J21"
$ws.Range("J21").Value = 11.282963378125267
Copy-Format "I8" "J22"
$ws.Range("J22").Value = 25.042808754677555
Copy-Format "I8" "J23"
$ws.Range("J23").Value = 3.2011163356916352
Copy-Format "I8" "J24"
$ws.Range("J24").Value = 13.523574517571838
Copy-Format "I8" "J25"
$ws.Range("J25").Value = 6.1196997869329204
Copy-Format "I8" "J26"
$ws.Range("J26").Value = 5.9488136666578013

# Bottom total row - already has the border + number format this column
# needs, so clone straight from the matching I27 cell.
Copy-Format "I27" "J27"
$ws.Range("J27").Value = 5.2451982064110645

$excel.CutCopyMode = $false

# Restore the selection that was active when the workbook was last saved.
$ws.Range("N8").Select() | Out-Null
